$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.277.31'
$ws.Range("E2").Value = '  +0.04%  '
$ws.Range("D3").Value = '1.864.00'
$ws.Range("E3").Value = '  -0.29%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.28'
$ws.Range("E5").Value = '  +0.58%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.05%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4709'
$ws.Range("E7").Value = '  +0.98%  '
$ws.Range("E8").Value = '  +2.40%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06540'
$ws.Range("E9").Value = '  -0.13%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.97'
$ws.Range("E10").Value = '  +3.87%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07941'
$ws.Range("E11").Value = '  +0.91%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '97.96'
$ws.Range("E12").Value = '  +0.09%  '
$ws.Range("D13").Value = '1.865.33'
$ws.Range("E13").Value = '  -0.26%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.149'
$ws.Range("E14").Value = '  +0.73%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6809'
$ws.Range("E15").Value = '  +1.11%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '262.93'
$ws.Range("E16").Value = '  -6.30%  '
$ws.Range("D17").Value = '30.265.35'
$ws.Range("E17").Value = '  +0.03%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.75'
$ws.Range("E18").Value = '  +8.50%  '
$ws.Range("E19").Value = '  +0.02%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007474'
$ws.Range("E20").Value = '  +2.64%  '
$ws.Range("D21").Value = '2.107.22'
$ws.Range("E21").Value = '  -0.25%  '
$ws.Range("E22").Value = '  -0.10%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.271'
$ws.Range("E23").Value = '  -4.38%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.180'
$ws.Range("E24").Value = '  +0.21%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '167.37'
$ws.Range("E25").Value = '  +1.60%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.196'
$ws.Range("E26").Value = '  -0.28%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.92'
$ws.Range("E27").Value = '  -1.46%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.950'
$ws.Range("E28").Value = '  +1.16%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.394'
$ws.Range("E29").Value = '  +1.44%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09936'
$ws.Range("E30").Value = '  +2.25%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.349'
$ws.Range("E31").Value = '  -1.69%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.470'
$ws.Range("E32").Value = '  -0.22%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.021'
$ws.Range("E33").Value = '  -2.17%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04716'
$ws.Range("E34").Value = '  +0.53%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.130'
$ws.Range("E35").Value = '  +0.76%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7000'
$ws.Range("E36").Value = '  -0.89%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.707'
$ws.Range("E37").Value = '  -0.83%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01881'
$ws.Range("E38").Value = '  +1.36%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.620'
$ws.Range("E39").Value = '  +3.12%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.320'
$ws.Range("E40").Value = '  +0.63%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '73.86'
$ws.Range("E41").Value = '  +0.46%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.942'
$ws.Range("E42").Value = '  -0.26%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8422'
$ws.Range("E43").Value = '  -0.44%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9997'
$ws.Range("E44").Value = '  -0.12%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4160'
$ws.Range("E45").Value = '  -0.25%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '103.26'
$ws.Range("E46").Value = '  -0.64%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.160'
$ws.Range("E47").Value = '  -0.44%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '949.59'
$ws.Range("E48").Value = '  +1.87%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.209'
$ws.Range("E49").Value = '  +0.53%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.17'
$ws.Range("E50").Value = '  +0.35%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05662'
$ws.Range("E51").Value = '  +0.59%  '
